$d = $word.ActiveDocument

$replacements = @(
    @{old="2024-03-21 Thursday"; new="2024-03-22 Friday"},
    @{old="582÷3="; new="885÷4="},
    @{old="598÷4="; new="669÷3="},
    @{old="299÷7="; new="246÷6="},
    @{old="593÷9="; new="796÷7="},
    @{old="289÷4="; new="719÷4="},
    @{old="944÷3="; new="580÷4="},
    @{old="114÷7="; new="978÷6="},
    @{old="780÷5="; new="291÷9="},
    @{old="115÷9="; new="776÷3="},
    @{old="696÷6="; new="821÷8="},
    @{old="167÷5="; new="697÷6="},
    @{old="572÷6="; new="724÷3="},
    @{old="813÷6="; new="816÷9="},
    @{old="670÷8="; new="671÷2="},
    @{old="150÷6="; new="851÷6="},
    @{old="396÷2="; new="180÷9="},
    @{old="741÷8="; new="585÷9="},
    @{old="207÷6="; new="442÷7="},
    @{old="850÷9="; new="671÷8="},
    @{old="917÷6="; new="134÷8="},
    @{old="411÷2="; new="748÷7="},
    @{old="491÷9="; new="956÷4="},
    @{old="635÷7="; new="544÷7="},
    @{old="629÷8="; new="638÷5="},
    @{old="522÷5="; new="715÷2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
